$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# --- Department column (C): replace the constant "SHELDON SCHOOL OF HOSPITALITY"
# with the course's actual department, grouped by course type ---
$ws.Range("C2:C5").Value = "Hospitality"
$ws.Range("C6:C9").Value = "Cookery"
$ws.Range("C10:C13").Value = "Patisserie and Baking"
$ws.Range("C14:C15").Value = "Event Management"
$ws.Range("C16:C17").Value = "Travel and Tourism"
$ws.Range("C18:C28").Value = "Packages"

# --- Location / locationDetail columns (M / N): split out the
# "(Currently not accepting enrolments)" qualifier into its own column ---
$ws.Range("M12").Value = "NSW/QLD/TAS"
$ws.Range("N12").Value = "Currently not accepting enrolments"

$ws.Range("M18").Value = "NSW/QLD/TAS"
$ws.Range("N18").Value = "Currently not accepting enrolments"

$ws.Range("M23").Value = "NSW/QLD/TAS"
$ws.Range("N23").Value = "Currently not accepting enrolments"

# --- Note column (R): the "Promotion valid until 31th Dec 2021" note has
# expired, clear it from every course row ---
$ws.Range("R2:R28").Value = ""
